# Update the analysis worksheet with new/refreshed model benchmark rows.
# Three previously-listed models (Qwen2, Phi-3, DeepSeek-V2) are replaced
# with their updated versions (Qwen2.5, Phi-3.5-MoE, DeepSeek-V2.5), and
# CodeQwen1.5 is replaced with Qwen2.5-Coder, along with refreshed scores.
# All dependent formulas (Max/Mean/Std columns, Improvement Analysis table,
# summary statistics, etc.) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Qwen2 -> Qwen2.5 ---
$ws.Range("C8").Value = "Qwen2.5"
$ws.Range("D8").Value = 75.923487071692406
$ws.Range("E8").Value = 86.040156752479803
$ws.Range("F8").Value = 85.821885287475496
$ws.Range("G8").Value = 78.710816662280493
$ws.Range("H8").Value = 86.956808148753694
$ws.Range("I8").Value = 85.546383356422396
$ws.Range("J8").Value = 83.1665895465174
$ws.Range("K8").Value = 4.5115418623013399

# --- Row 9: Phi-3 -> Phi-3.5-MoE ---
$ws.Range("C9").Value = "Phi–3.5-MoE"
$ws.Range("D9").Value = 50.553712395398797
$ws.Range("E9").Value = 74.284368466285997
$ws.Range("F9").Value = 74.622316617628201
$ws.Range("G9").Value = 71.019966600588802
$ws.Range("H9").Value = 75.447000426318297
$ws.Range("I9").Value = 75.674632955774698
$ws.Range("J9").Value = 70.266999576999197
$ws.Range("K9").Value = 8.9476588904403602

# --- Row 10: DeepSeek-V2 -> DeepSeek-V2.5 ---
$ws.Range("C10").Value = "DeepSeek–V2.5"
$ws.Range("D10").Value = 70.766842096210695
$ws.Range("E10").Value = 83.671431510681103
$ws.Range("F10").Value = 86.842415329459101
$ws.Range("G10").Value = 71.415797620888100
$ws.Range("H10").Value = 71.472025509189294
$ws.Range("I10").Value = 81.962624906808202
$ws.Range("J10").Value = 77.688522828872706
$ws.Range("K10").Value = 6.7839305479916101

# --- Row 13: CodeQwen1.5 -> Qwen2.5-Coder ---
$ws.Range("C13").Value = "Qwen2.5-Coder"
$ws.Range("D13").Value = 60.168481406952999
$ws.Range("E13").Value = 76.978833995825994
$ws.Range("F13").Value = 71.632966896838198
$ws.Range("G13").Value = 58.769922014162503
$ws.Range("H13").Value = 73.518812099485302
$ws.Range("I13").Value = 81.548498617449795
$ws.Range("J13").Value = 70.436252505119100
$ws.Range("K13").Value = 8.6400711500808693

# --- Mirror the renamed models in the "Improvement Analysis" table below ---
$ws.Range("C31").Value = "Qwen2.5"
$ws.Range("C32").Value = "Phi–3.5-MoE"
$ws.Range("C33").Value = "DeepSeek–V2.5"
$ws.Range("C36").Value = "Qwen2.5-Coder"

$excel.CalculateFull()

# Move the active selection, matching the cell the author last had selected.
[void]$ws.Range("J8").Select()
